$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "3HkSLidfgeLyM1izEZMvB4eKHi94U4HWbBdfvY48Vpq3"
$ws.Range("B4").Value = 0.2
$ws.Range("C4").Value = "5Lv2jkvax2bGZmPL7tsuDcdBJ72dtMSCaxXp1HVQkRrfDsyFqq9A98SmV9DsKm5m1f4kfWU6mKZTZPdHHnKzpG39"

$ws.Range("A5").Value = "3HkSLidfgeLyM1izEZMvB4eKHi94U4HWbBdfvY48Vpq3"
$ws.Range("B5").Value = 0.2
$ws.Range("C5").Value = "5mepBeyQa3hY5XvM2sWmdBkVCCEbeASmMbu95toqHM61Y22cFUxpksP1v8UTphorcG3vPtxgsX4JkB5rnLPox3rv"

$ws.Range("A6").Value = "3HkSLidfgeLyM1izEZMvB4eKHi94U4HWbBdfvY48Vpq3"
$ws.Range("B6").Value = 0.2
$ws.Range("C6").Value = "2Ng1fSefsqzd6ZwCwgmQU3yKQG7eVNWM68zDfomQV8Z2RWcP8kwp33gNXHfS9PrAbib136TAi7DBBwjXL3Bq3SXo"
